# Apply "tareas+ test de clase en verde" commit:
#  - Rename task in row 31 from "Crear nuevos tests" to
#    "Pasar Test pendiente de clase a verde" and update its hours
#    (estimated 3 -> 1, real "??" -> 1).
#  - Insert a new row 32 for task "Hacer test" (3 estimated hours,
#    author Adolfo, date 2016-10-11) that was split out of row 31.
#  - Move the active selection from D33 to A33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 31: task text + estimated/real hours.
$ws.Cells.Item(31, 1).Value = "Pasar Test pendiente de clase a verde"
$ws.Cells.Item(31, 2).Value = 1
$ws.Cells.Item(31, 3).Value = 1

# Make room for the new row 32 (split-off "Hacer test" task) while keeping
# the existing row 33 placeholder exactly where it is: inserting a row at
# 32 picks up the formatting of row 31 above it (matching the author's
# styles) but also pushes row 33 down to 34, so immediately delete the
# (now empty) row 33 that results to shift row 34 back up to 33.
$ws.Rows.Item(32).Insert()
$ws.Rows.Item(33).Delete()

$ws.Cells.Item(32, 1).Value = "Hacer test"
$ws.Cells.Item(32, 2).Value = 3
$ws.Cells.Item(32, 4).Value = "Adolfo"
$ws.Cells.Item(32, 5).Value = 42654

# Move the selection to A33 (was D33).
$ws.Range("A33").Select()
